$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new headers, copying the existing header formatting (bold,
# centered, bordered) from an existing header cell so the new cells share
# the same style as B1:H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (I0, IF) for rows 2-35.
$I = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,5,4,9,7,8,9,1,1,1,1,1,5,5,7,3)
$J = @(6,5,5,7,6,6,6,6,7,5,4,7,5,5,5,6,4,6,5,7,6,9,8,9,9,4,4,5,4,4,5,5,7,5)

for ($idx = 0; $idx -lt $I.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $I[$idx]
    $ws.Cells.Item($row, 10).Value = $J[$idx]
}
